$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing B2:G10 values down one row to B3:G11 (column A labels
# stay fixed in place), then write fresh values into B2:G2.
$ws.Range("B2:G10").Copy()
$ws.Range("B3:G11").PasteSpecial(-4163)  # xlPasteValues = -4163
$excel.CutCopyMode = $false

$ws.Range("B2").Value = 0.1783908196033299
$ws.Range("C2").Value = 0.3606156554386025
$ws.Range("D2").Value = 0.2599511937740667
$ws.Range("E2").Value = 0.5098540906711122
$ws.Range("F2").Value = 0.4943913024279584
$ws.Range("G2").Value = 15
